# Generate Report for Handoff
# The CI report workbook is regenerated on each run: rows whose status is
# "Ready for handoff" (plus the "Handback transform failed" row) get their
# latest-handoff timestamp refreshed to the time of this run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column (D) ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("D7").Value  = "2016-20-12 06:20:08"
$ws.Range("D10").Value = "2016-20-12 06:20:08"
$ws.Range("D11").Value = "2016-20-12 06:20:08"
$ws.Range("D12").Value = "2016-20-12 06:20:08"
$ws.Range("D13").Value = "2016-20-12 06:20:08"
$ws.Range("D14").Value = "2016-20-12 06:20:08"
$ws.Range("D15").Value = "2016-20-12 06:20:08"
$ws.Range("D16").Value = "2016-20-12 06:20:08"

# --- zh-cn sheet: "Latest Handoff Datetime" column (E) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("E7").Value  = "2016-03-12 06:20:02"
$ws.Range("E10").Value = "2016-03-12 06:20:02"
$ws.Range("E11").Value = "2016-03-12 06:20:02"
$ws.Range("E12").Value = "2016-03-12 06:20:02"
$ws.Range("E13").Value = "2016-03-12 06:20:02"
$ws.Range("E14").Value = "2016-03-12 06:20:02"
$ws.Range("E15").Value = "2016-03-12 06:20:02"
$ws.Range("E16").Value = "2016-03-12 06:20:02"

# --- de-de sheet: "Latest Handoff Datetime" column (E) ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("E7").Value  = "2016-03-12 06:20:08"
$ws.Range("E10").Value = "2016-03-12 06:20:08"
$ws.Range("E11").Value = "2016-03-12 06:20:08"
$ws.Range("E12").Value = "2016-03-12 06:20:08"
$ws.Range("E13").Value = "2016-03-12 06:20:08"
$ws.Range("E14").Value = "2016-03-12 06:20:08"
$ws.Range("E15").Value = "2016-03-12 06:20:08"
$ws.Range("E16").Value = "2016-03-12 06:20:08"
